$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("AG18").Value = 66.3
$ws.Range("AH18").Value = 192
